# Insert a new data row at row 132 (pushing the existing rows 132..209 down
# to 133..210) on the single worksheet of the workbook, and populate the new
# row with the Limache / Ají "Primera" record. This mirrors the diff, which
# shows every row from 132 through 209 shifting its data down by one row and
# a brand-new row 210 (formerly 209) appearing at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 132..209 down to 133..210, leaving row 132 blank and ready for
# the new record (matches Excel's normal "insert row" behaviour, including
# carrying the row's formatting down with it).
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the new record's values.
$ws.Range("A132").Value = 4
$ws.Range("B132").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C132").Value = "Los Lagos"
$ws.Range("D132").Value = 44596
$ws.Range("E132").Value = 10
$ws.Range("F132").Value = 100112021
$ws.Range("G132").Value = "Ají"
$ws.Range("H132").Value = "Inferno"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 150
$ws.Range("K132").Value = 25000
$ws.Range("L132").Value = 26000
$ws.Range("M132").Value = 25533
$ws.Range("N132").Value = "$/caja 15 kilos"
$ws.Range("O132").Value = "Limache"
$ws.Range("P132").Value = 1702
$ws.Range("Q132").Value = 15
$ws.Range("R132").Value = "Hortaliza"
